$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'258.68"
$ws.Range("D3").Value = "'26.86"
$ws.Range("E3").Value = "'-0.72%"
$ws.Range("E4").Value = "'1.57%"
$ws.Range("D5").Value = "'0.05986"
$ws.Range("E5").Value = "'1.30%"
$ws.Range("D6").Value = "'6.641"
$ws.Range("E6").Value = "'0.27%"
$ws.Range("D7").Value = "'0.8573"
$ws.Range("E7").Value = "'-0.32%"
$ws.Range("D8").Value = "'0.9234"
$ws.Range("E8").Value = "'-1.16%"
$ws.Range("D9").Value = "'0.1386"
$ws.Range("E9").Value = "'-1.50%"
$ws.Range("D10").Value = "'0.04539"
$ws.Range("E10").Value = "'25.06%"
$ws.Range("D11").Value = "'0.07008"
$ws.Range("E11").Value = "'-1.15%"
$ws.Range("D12").Value = "'0.03053"
$ws.Range("E12").Value = "'-5.53%"
$ws.Range("D13").Value = "'0.09116"
$ws.Range("E13").Value = "'-0.94%"
$ws.Range("D14").Value = "'0.001528"
$ws.Range("E14").Value = "'-1.13%"
$ws.Range("D15").Value = "'0.0006053"
$ws.Range("E15").Value = "'-94.17%"
$ws.Range("D16").Value = "'0.006130"
$ws.Range("E16").Value = "'0.32%"
$ws.Range("D17").Value = "'3.441"
$ws.Range("E17").Value = "'-2.11%"
$ws.Range("E18").Value = "'-1.48%"
$ws.Range("E19").Value = "'-2.16%"
$ws.Range("E20").Value = "'1.64%"
$ws.Range("E21").Value = "'0.86%"
$ws.Range("D22").Value = "'4.027"
$ws.Range("E22").Value = "'4.68%"
$ws.Range("D23").Value = "'0.04219"
$ws.Range("E23").Value = "'-0.05%"
$ws.Range("D25").Value = "'0.004030"
$ws.Range("D26").Value = "'0.0001199"
$ws.Range("E26").Value = "'-0.08%"
$ws.Range("E27").Value = "'-11.65%"
$ws.Range("D40").Value = "'0.03829"
$ws.Range("E40").Value = "'-0.02%"
$ws.Range("E41").Value = "'1.18%"
$ws.Range("D42").Value = "'0.003757"
$ws.Range("E42").Value = "'-40.09%"
$ws.Range("D43").Value = "'0.002419"
$ws.Range("E43").Value = "'9.92%"
$ws.Range("E44").Value = "'33.62%"
$ws.Range("E45").Value = "'-6.03%"
$ws.Range("E46").Value = "'-0.08%"
$ws.Range("E47").Value = "'-17.04%"
$ws.Range("D48").Value = "'0.1869"
$ws.Range("E48").Value = "'171.99%"
$ws.Range("E49").Value = "'-0.08%"
$ws.Range("E50").Value = "'-0.08%"
